$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 3
    4  = 3
    5  = 3
    6  = 3
    8  = 1
    9  = 6
    10 = -5
    11 = 5
    12 = -2
    13 = 11
    14 = -4
    16 = 3
    17 = 3
    18 = 1
    19 = -1
    20 = 2
    21 = -2
    22 = 7
    23 = -1
    24 = 2
    25 = 6
    26 = 1
    27 = -5
    28 = 2
    29 = -2
    30 = -2
    31 = -1
    33 = 4
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = -1
    39 = 5
    40 = -3
    41 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
